# Weekly refresh: insert the newest daily record at the top of the data
# block (row 281) and push every existing record down by one row. The
# oldest record that falls off the bottom of the original range becomes
# the new last row (340).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 281..339 down to 282..340 (also extends the used range /
# dimension to row 340 and carries the date-formatted style on column D).
$ws.Rows.Item(281).Insert()

# Populate the newly inserted row with this week's record.
$ws.Range("A281").Value = 9
$ws.Range("B281").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C281").Value = "Metropolitana"
$ws.Range("D281").Value = 44641
$ws.Range("E281").Value = 13
$ws.Range("F281").Value = 100112044
$ws.Range("G281").Value = "Perejil"
$ws.Range("H281").Value = "Sin especificar"
$ws.Range("I281").Value = "Primera"
$ws.Range("J281").Value = 43
$ws.Range("K281").Value = 12000
$ws.Range("L281").Value = 14000
$ws.Range("M281").Value = 13023
$ws.Range("N281").Value = "$/docena de atados"
$ws.Range("O281").Value = "Región Metropolitana"
$ws.Range("P281").Value = 4341
$ws.Range("Q281").Value = 3
$ws.Range("R281").Value = "Hortaliza"
